# elapsed time y cpu
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy formatting (bold, border, center alignment) from the existing
# header cell F1 onto the two new header cells G1:H1
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Elapsed Time / CPU values for each data row
$elapsed = 0.4794827245333484
$cpu = 0.996

$ws.Range("G2").Value = $elapsed
$ws.Range("H2").Value = $cpu

$ws.Range("G3").Value = $elapsed
$ws.Range("H3").Value = $cpu

$ws.Range("G4").Value = $elapsed
$ws.Range("H4").Value = $cpu

$ws.Range("G5").Value = $elapsed
$ws.Range("H5").Value = $cpu

# Minor re-computed precision refresh on existing MSE/MAE columns
# (negligible last-digit differences from the re-run pipeline)
$ws.Range("B3").Value = 0.07514644587374571
$ws.Range("D3").Value = 0.2119198634755612

$ws.Range("B4").Value = 0.04215534119371416
$ws.Range("D4").Value = 0.1361288253571671

$ws.Range("B5").Value = 0.07796894984218657
$ws.Range("D5").Value = 0.1911874935925048
